# Seguridad con token jwt
# Adds a new "company" worksheet (code/name lookup table) right before the
# existing "mandarEncuesta" sheet, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

$mandar = $wb.Worksheets.Item("mandarEncuesta")

# Create the new sheet unpositioned first -- inserting it directly at a
# specific position confuses later Copy/PasteSpecial calls against it, so
# build it in place and reposition with Move() once it is fully populated.
$company = $wb.Worksheets.Add()
$company.Name = "company"

# Header row.
$company.Range("A1").Value = "code"
$company.Range("B1").Value = "name"

# Data row.
$company.Range("A2").Value = "E0001"
$company.Range("B2").Value = "Sun Sol"

# Match the header formatting already used across the workbook (white bold
# text on a dark grey fill) by copying it from an existing header cell
# instead of inventing a brand-new style.
$mandar.Range("A1:B1").Copy()
$company.Range("A1:B1").PasteSpecial(-4122)
$company.Range("A1").Select() | Out-Null

# Move the new sheet right after "Encuestas" (i.e. right before
# "mandarEncuesta"), matching the target sheet order:
#   divisiónTerritorial, divisiónServicios, Usuarios, Encuestas, company, mandarEncuesta
$encuestas = $wb.Worksheets.Item("Encuestas")
$company.Move($null, $encuestas)

# Re-fetch sheet references by name after the move/reorder.
$company = $wb.Worksheets.Item("company")
$mandar = $wb.Worksheets.Item("mandarEncuesta")

# Column widths (A:B wide for code/name, rest narrow) matching the new
# sheet's layout.
$company.Range("A1:B1").ColumnWidth = 21.85
$company.Range("C1").EntireColumn.ColumnWidth = 7.35

# Make "company" the active/selected tab with C8 highlighted, and restore
# "mandarEncuesta" selection back to A1 (it is no longer the active tab).
$mandar.Range("A1").Select() | Out-Null
$company.Activate() | Out-Null
$company.Range("C8").Select() | Out-Null
